# Adicionados balancos concatenados em uma unica planilha.
#
# A new quarter (31/12/2008, already present as the header of column AI) was
# missing its own flow figures for the income-statement rows (58-79): AI
# actually held the figure that belonged to the following quarter (31/03/2009,
# column AJ), and so on down the line. This inserts the correct AI-quarter
# value and shifts the old AI..AQ values one column to the right (AJ..AR) -
# i.e. a "shift cells right" insert scoped to a single row - and does the
# same for the CH..CP block (31/12/2021 onward) where that block has data.
#
# NOTE: this sandboxed engine's Range.Insert() shifts the whole column (every
# row), not just the target row, so a plain Range.Insert(xlShiftToRight) call
# is unusable here. Instead each row is shifted manually, cell by cell, only
# within that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlShiftToRight = -4161  # kept only for documentation of the semantics being emulated

# Column numbers (1-based) for the two affected blocks.
$AI = 35   # column AI
$AQ = 43   # column AQ
$CH = 86   # column CH
$CP = 94   # column CP

function Insert-ShiftRight {
    # Emulates Range(<startCol>:<endCol> on row $row).Insert(xlShiftToRight)
    # by copying values one cell to the right, starting from the rightmost
    # column so nothing is overwritten before it is read, then writes
    # $newValue into the freed-up first cell. Scoped to a single row only.
    param($row, $startCol, $endCol, $newValue)

    for ($c = $endCol; $c -ge $startCol; $c--) {
        $v = $ws.Cells.Item($row, $c).Value()
        $ws.Cells.Item($row, $c + 1).Value = $v
    }
    $ws.Cells.Item($row, $startCol).Value = $newValue
}

# New value to insert at AI<row> for each affected income-statement row
# (old AI..AQ shift right to AJ..AR).
$newAI = [ordered]@{
    58 = 170133.984
    59 = -107612
    60 = 62522
    61 = -29843
    62 = 2997
    63 = -12135
    64 = -16650
    65 = -313
    66 = 1432
    67 = -5174
    69 = 32679
    70 = -233
    71 = -1037
    72 = 755
    73 = 32446
    74 = 2835
    75 = -24502
    76 = -317
    79 = 10462
}

# New value to insert at CH<row> for each affected income-statement row
# (old CH..CP shift right to CI..CQ). Only rows whose CH..CP block actually
# carried data get an entry here.
$newCH = [ordered]@{
    58 = 169527.008
    59 = -68483
    60 = 101044
    62 = 16018
    63 = -37352
    64 = -19443
    65 = -6138
    66 = 799
    67 = -17239
    68 = 2388
    73 = 40077
    74 = -15594
    76 = -1580
    79 = 22903
}

foreach ($row in 58..79) {
    if ($newAI.Contains($row)) {
        Insert-ShiftRight $row $AI $AQ $newAI[$row]
    }
    if ($newCH.Contains($row)) {
        Insert-ShiftRight $row $CH $CP $newCH[$row]
    }
}
